$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise previously-published monthly M2/FX figures (rows 122-218) ---
# Each (B = local-currency M2, D = USD-converted M2) pair is updated per the
# refreshed TradingView feed snapshot; C (FX rate) is unchanged.
$ws.Cells.Item(122,2).Value = 261572700000
$ws.Cells.Item(122,4).Value = 70634235255.99483
$ws.Cells.Item(123,2).Value = 259192800000
$ws.Cells.Item(123,4).Value = 65974189935.6021
$ws.Cells.Item(124,2).Value = 258559200000
$ws.Cells.Item(124,4).Value = 65332322619.76955
$ws.Cells.Item(125,2).Value = 255274400000
$ws.Cells.Item(125,4).Value = 62184697084.11488
$ws.Cells.Item(131,2).Value = 263647500000
$ws.Cells.Item(131,4).Value = 66715800394.75682
$ws.Cells.Item(134,2).Value = 286126300000
$ws.Cells.Item(134,4).Value = 69012614568.25856
$ws.Cells.Item(135,2).Value = 283933800000
$ws.Cells.Item(135,4).Value = 67701614249.26678
$ws.Cells.Item(136,2).Value = 283623100000
$ws.Cells.Item(136,4).Value = 69108942495.12671
$ws.Cells.Item(137,2).Value = 280655000000
$ws.Cells.Item(137,4).Value = 71570102514.40811
$ws.Cells.Item(142,2).Value = 295604600000
$ws.Cells.Item(142,4).Value = 74183045573.17809
$ws.Cells.Item(143,2).Value = 295922500000
$ws.Cells.Item(143,4).Value = 74799681512.56256
$ws.Cells.Item(144,2).Value = 296506400000
$ws.Cells.Item(144,4).Value = 72371588967.53722
$ws.Cells.Item(146,2).Value = 314026000000
$ws.Cells.Item(146,4).Value = 72895378258.54825
$ws.Cells.Item(147,2).Value = 309824800000
$ws.Cells.Item(147,4).Value = 74464585286.12973
$ws.Cells.Item(148,2).Value = 311655400000
$ws.Cells.Item(148,4).Value = 72999180193.47433
$ws.Cells.Item(149,2).Value = 314511600000
$ws.Cells.Item(149,4).Value = 73714808043.87569
$ws.Cells.Item(150,2).Value = 318498700000
$ws.Cells.Item(150,4).Value = 76541948042.58488
$ws.Cells.Item(151,2).Value = 320618400000
$ws.Cells.Item(151,4).Value = 79071322876.59071
$ws.Cells.Item(152,2).Value = 321076200000
$ws.Cells.Item(152,4).Value = 80883766626.36035
$ws.Cells.Item(153,2).Value = 324218500000
$ws.Cells.Item(153,4).Value = 84278268780.86821
$ws.Cells.Item(154,2).Value = 329607000000
$ws.Cells.Item(154,4).Value = 85630001039.17696
$ws.Cells.Item(155,2).Value = 332065200000
$ws.Cells.Item(155,4).Value = 85414306659.46446
$ws.Cells.Item(156,2).Value = 336385900000
$ws.Cells.Item(156,4).Value = 85132968896.31261
$ws.Cells.Item(157,2).Value = 339280200000
$ws.Cells.Item(157,4).Value = 87126730181.55672
$ws.Cells.Item(158,2).Value = 350004800000
$ws.Cells.Item(158,4).Value = 90272567832.45641
$ws.Cells.Item(159,2).Value = 348823700000
$ws.Cells.Item(159,4).Value = 93148819696.64603
$ws.Cells.Item(160,2).Value = 352411100000
$ws.Cells.Item(160,4).Value = 92319466638.72371
$ws.Cells.Item(161,2).Value = 351238500000
$ws.Cells.Item(161,4).Value = 93070431119.00156
$ws.Cells.Item(162,2).Value = 354917800000
$ws.Cells.Item(162,4).Value = 91978593826.93654
$ws.Cells.Item(163,2).Value = 356539500000
$ws.Cells.Item(163,4).Value = 89526553672.31639
$ws.Cells.Item(164,2).Value = 362385300000
$ws.Cells.Item(164,4).Value = 90843874558.17102
$ws.Cells.Item(165,2).Value = 360221900000
$ws.Cells.Item(165,4).Value = 91149266194.33199
$ws.Cells.Item(170,2).Value = 381075300000
$ws.Cells.Item(170,4).Value = 93916428430.59937
$ws.Cells.Item(171,2).Value = 382602200000
$ws.Cells.Item(171,4).Value = 92581474132.50739
$ws.Cells.Item(172,2).Value = 384958100000
$ws.Cells.Item(172,4).Value = 92287320499.60445
$ws.Cells.Item(173,2).Value = 383090000000
$ws.Cells.Item(173,4).Value = 90073123133.71423
$ws.Cells.Item(182,2).Value = 422631600000
$ws.Cells.Item(182,4).Value = 99109255915.39055
$ws.Cells.Item(192,2).Value = 469280100000
$ws.Cells.Item(192,4).Value = 112585792428.3864
$ws.Cells.Item(194,2).Value = 487349900000
$ws.Cells.Item(194,4).Value = 123030874482.4801
$ws.Cells.Item(195,2).Value = 490302300000
$ws.Cells.Item(195,4).Value = 122236369075.8146
$ws.Cells.Item(196,2).Value = 496963100000
$ws.Cells.Item(196,4).Value = 123215010041.4053
$ws.Cells.Item(197,2).Value = 499199700000
$ws.Cells.Item(197,4).Value = 119041302015.0232
$ws.Cells.Item(206,2).Value = 564423000000
$ws.Cells.Item(206,4).Value = 129889768490.8179
$ws.Cells.Item(210,2).Value = 569711700000
$ws.Cells.Item(210,4).Value = 121551461489.2255
$ws.Cells.Item(212,2).Value = 569309400000
$ws.Cells.Item(212,4).Value = 120803233815.0105
$ws.Cells.Item(215,2).Value = 581768700000
$ws.Cells.Item(215,4).Value = 115151557737.2234
$ws.Cells.Item(218,2).Value = 603199600000
$ws.Cells.Item(218,4).Value = 131769141708.7184

# --- Append three new monthly rows (221:223) for Mar/Apr/May 2023 ---
# Copy the date-cell format (border/bold/center/date-time numfmt) from the
# last existing row (A220) down into the new date cells so they match the
# sheet's existing style (s="2") instead of minting a new style index.
$ws.Range("A220").Copy() | Out-Null
$ws.Range("A221:A223").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(221,1).Value = 44986
$ws.Cells.Item(221,2).Value = 613926400000
$ws.Cells.Item(221,3).Value = 0.2205168915938961
$ws.Cells.Item(221,4).Value = 135381141395.4309

$ws.Cells.Item(222,1).Value = 45017
$ws.Cells.Item(222,2).Value = 618680400000
$ws.Cells.Item(222,3).Value = 0.2248302531588651
$ws.Cells.Item(222,4).Value = 139098070956.4279

$ws.Cells.Item(223,1).Value = 45047
$ws.Cells.Item(223,2).Value = 624790700000
$ws.Cells.Item(223,3).Value = 0.2154151049071561
$ws.Cells.Item(223,4).Value = 134589354185.5155
